# feature | i-475 | Ajustes importacion masiva kits
# Adds "Descripcion" and "Nombre secundario" columns (J, K) to the item-sets
# import template, with sample data for the two existing kit rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells, then fill in column J (Descripcion) and
# column K (Nombre secundario) sample data, matching the order the strings
# were authored in (shared-string table order).
$ws.Range("J1").Value = "Descripcion"
$ws.Range("K1").Value = "Nombre secundario"

$ws.Range("J2").Value = "desc 1"
$ws.Range("J3").Value = "desc 2"

$ws.Range("K2").Value = "nombre sec 1"
$ws.Range("K3").Value = "nombre sec 2"

# Widen column K so the new "Nombre secundario" header/content fits
$ws.Columns.Item(11).ColumnWidth = 18.5703125

# Leave the selection on the first new cell, matching the saved view state
$ws.Range("J1").Select()
